$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 - set the value then copy the formatting from the
# neighbouring header cell (G1) so it picks up the same bold/bordered
# header style used by the rest of row 1.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New "Save" data column, H2:H8
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
